$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data changes -----------------------------------------------------

# Row 85: add a note in column B explaining the trial observable
$ws.Range("B85").Value = "This is a trial obs, which is the per-person version of wage_rgd_demean_obs. I compare these two variables"

# Row 79: fix casing of the id (I_rgpc_obs -> i_rgpc_obs)
$ws.Range("A79").Value = "i_rgpc_obs"

# Row 86 (new): government debt observable
$ws.Range("A86").Value = "govdebt_rcpc_obs"
$ws.Range("B86").Value = "Constructed change of log real per capita government debt series,  demeaned"
$ws.Range("C86").Value = "government debt, constructed"
$ws.Range("D86").Value = "RENTIN-CPROFIT-W255RC1Q027SBEA-PROPINC-A074RC1Q027SBEA-W071RC1Q027SBEA-WASCUR-PROPINC-COE-W780RC1Q027SBEA-B249RC1Q027SBEA-B075RC1Q027SBEA-GDPCTPI-CNP16OV-A957RC1Q027SBEA-A787RC1Q027SBEA-AD08RC1Q027SBEA-A918RC1Q027SBEA-MVGFD027MNFRBDAL-W014RC1Q027SBEA-W011RC1Q027SBEA-W020RC1Q027SBEA-B232RC1Q027SBEA-B096RC1Q027SBEA-W006RC1Q027SBEA-W780RC1Q027SBEA-W009RC1Q027SBEA-B097RC1Q027SBEA-A091RC1Q027SBEA"

# --- Column widths (B and C were narrowed) -----------------------------
# The saved OOXML column width is ColumnWidth + 5/6, snapped to the nearest
# 1/6 of a character; choose inputs landing on the closest achievable step
# to the target widths 46.88671875 / 37.109375.
$ws.Columns("B").ColumnWidth = 46
$ws.Columns("C").ColumnWidth = 36.333333333333336

# --- Row heights (wrap-text rows got taller once columns got narrower) -
$ws.Rows("22").RowHeight = 57.6
$ws.Rows("23").RowHeight = 43.2
$ws.Rows("26").RowHeight = 28.8

# --- View state (best effort) ------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 73
$win.ScrollColumn = 3
$win.Zoom = 70
$win.Left = 3288
$win.Top = 2988
$win.Width = 16920
$win.Height = 9072
$ws.Range("D86").Select()
